$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.396.59"
$ws.Range("E2").Value = "  +0.29%  "
Set-TextValue "D3" "1.693.31"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.31%  "
Set-TextValue "D5" "219.04"
$ws.Range("E5").Value = "  +0.24%  "
Set-TextValue "D6" "0.5495"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("E7").Value = "  +0.28%  "
Set-TextValue "D8" "0.2746"
$ws.Range("E8").Value = "  +1.37%  "
Set-TextValue "D9" "0.06458"
$ws.Range("E9").Value = "  +0.44%  "
Set-TextValue "D11" "0.07679"
$ws.Range("E11").Value = "  +2.46%  "
Set-TextValue "D12" "1.704.99"
$ws.Range("E12").Value = "  +0.42%  "
Set-TextValue "D13" "4.544"
$ws.Range("E13").Value = "  -0.62%  "
Set-TextValue "D14" "0.5835"
$ws.Range("E14").Value = "  -0.40%  "
Set-TextValue "D15" "0.000008360"
$ws.Range("E15").Value = "  -1.85%  "
Set-TextValue "D16" "65.49"
$ws.Range("E16").Value = "  +1.48%  "
Set-TextValue "D17" "26.440.88"
$ws.Range("E17").Value = "  +0.23%  "
Set-TextValue "D18" "4.935"
$ws.Range("E18").Value = "  -0.26%  "
Set-TextValue "D19" "1.011"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  +0.74%  "
Set-TextValue "D21" "191.89"
$ws.Range("E21").Value = "  +1.17%  "
Set-TextValue "D22" "6.252"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +0.26%  "
Set-TextValue "D24" "148.86"
$ws.Range("E24").Value = "  +2.80%  "
Set-TextValue "D25" "0.1328"
$ws.Range("E25").Value = "  +7.65%  "
Set-TextValue "D26" "7.914"
$ws.Range("E26").Value = "  +2.67%  "
Set-TextValue "D27" "15.78"
$ws.Range("E27").Value = "  -0.54%  "
Set-TextValue "D28" "0.06288"
$ws.Range("E28").Value = "  -6.55%  "
Set-TextValue "D29" "1.383"
$ws.Range("E29").Value = "  +1.90%  "
Set-TextValue "D30" "1.333"
$ws.Range("E30").Value = "  +0.22%  "
Set-TextValue "D31" "3.605"
$ws.Range("E31").Value = "  +0.44%  "
Set-TextValue "D32" "3.612"
$ws.Range("E32").Value = "  +1.03%  "
Set-TextValue "D33" "1.684"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  +1.11%  "
Set-TextValue "D35" "0.6155"
$ws.Range("E35").Value = "  -1.43%  "
Set-TextValue "D36" "2.413"
$ws.Range("E36").Value = "  +0.70%  "
Set-TextValue "D37" "2.708"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  -3.01%  "
Set-TextValue "D40" "1.118.32"
$ws.Range("E40").Value = "  +0.20%  "
Set-TextValue "D41" "0.8848"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  -0.11%  "
Set-TextValue "D43" "101.57"
$ws.Range("E43").Value = "  +0.70%  "
Set-TextValue "D44" "1.844.27"
$ws.Range("E44").Value = "  +0.27%  "
Set-TextValue "D45" "57.52"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  -6.64%  "
Set-TextValue "D47" "1.012"
$ws.Range("E47").Value = "  +0.49%  "
Set-TextValue "D48" "8.190"
$ws.Range("E48").Value = "  -0.20%  "
Set-TextValue "D49" "0.05277"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D50" "6.117"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D51" "0.4305"
$ws.Range("E51").Value = "  -0.01%  "
